$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all touched cells to Text format to preserve exact string representation
# (prevents Excel auto-converting numeric-looking strings like "0.7560" into numbers)
$cells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "E31", "D32", "E32", "E33", "D34", "E34", "D35", "E35", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "E48", "B49", "C49", "D49", "E49", "B50", "C50", "D50", "E50", "D51", "E51")
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "26.895.65"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.812.85"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "309.15"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.4665"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("D8").Value = "0.3668"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "0.07351"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "0.8673"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "1.822.43"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "5.378"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "0.07084"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "91.65"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "1.003"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "0.000008696"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "1.001"
$ws.Range("D20").Value = "14.65"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "26.925.34"
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").Value = "5.298"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "10.62"
$ws.Range("E23").Value = "  -0.38%  "
$ws.Range("D24").Value = "2.044.80"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "1.893"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Value = "150.05"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").Value = "2.162"
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").Value = "18.27"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").Value = "5.275"
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("D30").Value = "115.58"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").Value = "0.7560"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").Value = "4.487"
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("D35").Value = "2.913"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "1.085"
$ws.Range("E37").Value = "  -1.98%  "
$ws.Range("D38").Value = "0.05282"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("D39").Value = "2.991"
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("D40").Value = "0.01949"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "7.214"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").Value = "0.5301"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").Value = "2.287"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("D44").Value = "0.1653"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "8.398"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("D46").Value = "0.4870"
$ws.Range("E46").Value = "  -2.66%  "
$ws.Range("D47").Value = "10.45"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "103.15"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.659"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").Value = "0.06288"
$ws.Range("E51").Value = "  -0.02%  "
